$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Categoria" column (D) was removed entirely; everything to its right
# shifts one column to the left (cell values, styles, column widths and
# data-validation ranges all move automatically with it).
$ws.Range("D1").EntireColumn.Delete() | Out-Null

# Rename the "limite_de_credito" header (now in column V) to "Limite de Credito".
$ws.Range("V1").Value = "Limite de Credito"

# Reflect the post-edit selection: whole column D (formerly column E), which is
# what's left standing where the user's cursor ended up after the deletion.
$ws.Range("D1:D1048576").Select() | Out-Null
